# NMRA_Region_Division_Map.xlsx — correct the "Lakeshores" division's
# region/division codes and drop the duplicate/obsolete "Former Lakeshores"
# row that it was mistakenly split off from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stray duplicate row (region 22 / division 3 / "Former Lakeshores")
# is removed first, while row 17 ("Lakeshores", still region 21 / division 15
# at this point) keeps the data range in its original, already-sorted order.
$ws.Rows.Item(21).Delete()

# Re-record the sort state over the new (now one-row-shorter) range so the
# saved sortState/sortCondition refs reflect the current extent.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3:A198"))
$ws.Sort.SortFields.Add($ws.Range("B3:B198"))
$ws.Sort.SetRange($ws.Range("A3:D198"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Now correct the "Lakeshores" row: it actually belongs to region 22
# (NFR), division 3 — not region 21 (NER), division 15.
$ws.Range("A17").Value = 22
$ws.Range("B17").Value = 3

# Match the author's final selection/view state.
$ws.Range("C20").Select()
